$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted before the existing row 54, which
# pushes the existing row 54..167 down to 55..168 (dimension grows from
# A1:R167 to A1:R168).
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new record's data.
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(54, 3).Value = 'La Araucanía'
$ws.Cells.Item(54, 4).Value = 45070
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = 100112035
$ws.Cells.Item(54, 7).Value = 'Bruselas (repollito)'
$ws.Cells.Item(54, 8).Value = 'Sin especificar'
$ws.Cells.Item(54, 9).Value = 'Primera'
$ws.Cells.Item(54, 10).Value = 20
$ws.Cells.Item(54, 11).Value = 28000
$ws.Cells.Item(54, 12).Value = 28000
$ws.Cells.Item(54, 13).Value = 28000
$ws.Cells.Item(54, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(54, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(54, 16).Value = 1867
$ws.Cells.Item(54, 17).Value = 15
$ws.Cells.Item(54, 18).Value = 'Hortaliza'
